$d = $word.ActiveDocument

# 1. Remove the page-break paragraph that sits between the "Chess" title
#    paragraph and the "0. Prelude" heading paragraph.
$pageBreakPara = $d.Paragraphs.Item(2)
$pageBreakPara.Range.Delete()

# 2. Append " Alpha" to the title paragraph ("Chess" -> "Chess Alpha").
#    Insert a temporary marker character after the new text so that the
#    insertion point used for the bookmark below is not sitting exactly on
#    a paragraph boundary (collapsed ranges right at a paragraph boundary
#    get attributed to the following paragraph instead of the current one).
$titlePara = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($titlePara.Range.End - 1, $titlePara.Range.End - 1)
$insertPoint.InsertAfter(" AlphaX")

# 3. Move the "_GoBack" bookmark from the end of the "5. Components and
#    implementation" heading to the end of the title paragraph (right after
#    the newly added " Alpha" text). Adding a bookmark with a name that
#    already exists relocates it instead of creating a duplicate.
$titlePara = $d.Paragraphs.Item(1)
$bmPos = $titlePara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary marker character ("X") now that the bookmark is
# safely anchored just before it.
$titlePara = $d.Paragraphs.Item(1)
$xRange = $d.Range($titlePara.Range.End - 2, $titlePara.Range.End - 1)
$xRange.Delete()

# 4. Refresh the "0. Prelude" run so the stale w:lastRenderedPageBreak
#    marker (left over from when this heading used to immediately follow
#    the page break) is dropped.
$d.Content.Find.Execute("0. Prelude", $false, $false, $false, $false, $false, $true, 1, $false, "0. Prelude", 2) | Out-Null

Write-Host "Title paragraph text:" $d.Paragraphs.Item(1).Range.Text
Write-Host "Paragraph count:" $d.Paragraphs.Count
